# Auto-generated Excel COM-interop script applying the diff to Sheet1.
# This mirrors data corrections made by an automated GH Action data refresh
# (updated vaccination / case counts for several Italian regions).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("L2").Value = 707439
$ws.Range("S2").Value = 35372
$ws.Range("T2").Value = 613189
$ws.Range("U2").Value = 306.59
$ws.Range("V2").Value = 919.78
$ws.Range("W2").Value = 1532.97
# Row 3
$ws.Range("K3").Value = 3647
$ws.Range("L3").Value = 270950
$ws.Range("S3").Value = 13548
$ws.Range("T3").Value = 290177
$ws.Range("U3").Value = 145.09
$ws.Range("V3").Value = 435.27
$ws.Range("W3").Value = 725.4400000000001
# Row 4
$ws.Range("K4").Value = 10194
$ws.Range("L4").Value = 942335
$ws.Range("S4").Value = 47117
$ws.Range("T4").Value = 982510
$ws.Range("U4").Value = 491.26
$ws.Range("V4").Value = 1473.77
$ws.Range("W4").Value = 2456.28
# Row 5
$ws.Range("K5").Value = 20190
$ws.Range("L5").Value = 3213947
$ws.Range("S5").Value = 160697
$ws.Range("T5").Value = 2626509
$ws.Range("U5").Value = 1313.25
$ws.Range("V5").Value = 3939.76
$ws.Range("W5").Value = 6566.27
# Row 6
$ws.Range("K6").Value = 24085
$ws.Range("L6").Value = 2370181
$ws.Range("S6").Value = 118509
$ws.Range("T6").Value = 2193877
$ws.Range("U6").Value = 1096.94
$ws.Range("V6").Value = 3290.82
$ws.Range("W6").Value = 5484.69
# Row 7
$ws.Range("K7").Value = 6029
$ws.Range("L7").Value = 656126
$ws.Range("T7").Value = 575433
$ws.Range("V7").Value = 863.15
$ws.Range("W7").Value = 1438.58
# Row 8
$ws.Range("L8").Value = 3394806
$ws.Range("T8").Value = 2495730
$ws.Range("U8").Value = 1247.87
# Row 9
$ws.Range("L9").Value = 783938
$ws.Range("T9").Value = 765064
$ws.Range("V9").Value = 1147.6
$ws.Range("W9").Value = 1912.66
# Row 10
$ws.Range("K10").Value = 35998
$ws.Range("L10").Value = 5722097
$ws.Range("S10").Value = 286105
$ws.Range("T10").Value = 4531000
$ws.Range("U10").Value = 2265.5
$ws.Range("V10").Value = 6796.5
$ws.Range("W10").Value = 11327.5
# Row 11
$ws.Range("K11").Value = 8163
$ws.Range("L11").Value = 833974
$ws.Range("S11").Value = 41699
$ws.Range("T11").Value = 709131
$ws.Range("U11").Value = 354.57
$ws.Range("V11").Value = 1063.7
$ws.Range("W11").Value = 1772.83
# Row 13
$ws.Range("L13").Value = 266717
$ws.Range("T13").Value = 280334
# Row 15
$ws.Range("K15").Value = 24357
$ws.Range("L15").Value = 2296035
$ws.Range("S15").Value = 114802
$ws.Range("T15").Value = 2091977
$ws.Range("U15").Value = 1045.99
$ws.Range("V15").Value = 3137.97
$ws.Range("W15").Value = 5229.94
# Row 16
$ws.Range("K16").Value = 11177
$ws.Range("L16").Value = 2240349
$ws.Range("S16").Value = 112017
$ws.Range("T16").Value = 1798599
$ws.Range("U16").Value = 899.3
$ws.Range("V16").Value = 2697.9
$ws.Range("W16").Value = 4496.5
# Row 17
$ws.Range("K17").Value = 5123
$ws.Range("L17").Value = 854953
$ws.Range("T17").Value = 786020
# Row 18
$ws.Range("K18").Value = 18327
$ws.Range("L18").Value = 2416021
$ws.Range("S18").Value = 120801
$ws.Range("T18").Value = 2545656
$ws.Range("U18").Value = 1272.83
$ws.Range("V18").Value = 3818.48
$ws.Range("W18").Value = 6364.14
# Row 19
$ws.Range("K19").Value = 23661
$ws.Range("L19").Value = 1852623
$ws.Range("S19").Value = 92631
$ws.Range("T19").Value = 1908341
$ws.Range("U19").Value = 954.17
$ws.Range("V19").Value = 2862.51
$ws.Range("W19").Value = 4770.85
# Row 20
$ws.Range("K20").Value = 2849
$ws.Range("L20").Value = 464101
$ws.Range("S20").Value = 23205
$ws.Range("T20").Value = 424117
$ws.Range("U20").Value = 212.06
$ws.Range("V20").Value = 636.1799999999999
$ws.Range("W20").Value = 1060.29
# Row 22
$ws.Range("K22").Value = 17910
$ws.Range("L22").Value = 2630769
$ws.Range("S22").Value = 131538
$ws.Range("T22").Value = 2353222
$ws.Range("U22").Value = 1176.61
$ws.Range("V22").Value = 3529.83
$ws.Range("W22").Value = 5883.05
$ws.Range("AD22").Value = 516
